$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits (tasks list was reworked) ---
$ws.Range("B5").Value = "Options menu (Sliders, toggles)"
$ws.Range("C5").Value = "Chatroom class"
$ws.Range("B6").Value = "Scores(Minus score if leave chatrooms be if red after sometime, touch wrong door)"
$ws.Range("C6").Value = "Enhanced touch detection for rooms, only red rooms can be touched to score"
$ws.Range("B7").Value = " touch with rooms"
$ws.Range("B8").Value = "Layout of pages"
$ws.Range("B9").Value = "Reset doors back to white if left untouched"

# Rows 10-11 no longer hold the last two tasks; clear them out but keep the rows.
$ws.Range("A10:H11").ClearContents()
$ws.Range("B10").Clear()

# --- Column width tweak (C got a bit wider) ---
$ws.Columns("C").ColumnWidth = 70.28515625

# --- Selection moved ---
$ws.Range("C10").Select()

$wb.Windows.Item(1).WindowState = -4143
